$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H2").Value = 848.25
$ws.Range("J2").Value = 798
$ws.Range("L2").Value = 798
$ws.Range("N2").Value = -1024
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H106").Value = 500000740
$ws.Range("I106").Value = 500000740
$ws.Range("K106").Value = 500000740
$ws.Range("M106").Value = -500000109
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H138").Value = 2851.1304
$ws.Range("J138").Value = 3241.7896
$ws.Range("L138").Value = 9725.3688
$ws.Range("N138").Value = -20005.3688

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -134
$ws.Range("H21").Value = 2071
$ws.Range("I21").Value = 1142
$ws.Range("K21").Value = 1142
$ws.Range("M21").Value = -768
$ws.Range("H30").Value = 1539.5
$ws.Range("I30").Value = 1539.5
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1539.5
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1389.5
$ws.Range("N30").ClearContents()
$ws.Range("H32").Value = 550.5454999999999
$ws.Range("I32").Value = 550.5454999999999
$ws.Range("K32").Value = 550.5454999999999
$ws.Range("M32").Value = -263.5454999999999
$ws.Range("H119").Value = 56961
$ws.Range("J119").Value = 56961
$ws.Range("L119").Value = 56961
$ws.Range("N119").Value = -66637
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 1628.4
$ws.Range("I86").Value = 1417.4615
$ws.Range("J86").Value = 2999.5
$ws.Range("K86").Value = 1417.4615
$ws.Range("L86").Value = 2999.5
$ws.Range("M86").Value = -294.4614999999999
$ws.Range("N86").Value = -5245.5
$ws.Range("H89").Value = 1628.4
$ws.Range("I89").Value = 1417.4615
$ws.Range("J89").Value = 2999.5
$ws.Range("K89").Value = 7087.307499999999
$ws.Range("L89").Value = 14997.5
$ws.Range("M89").Value = -1471.307499999999
$ws.Range("N89").Value = -26229.5
$ws.Range("H99").Value = 4838.9473
$ws.Range("I99").Value = 4424.857
$ws.Range("K99").Value = 4424.857
$ws.Range("M99").Value = -2926.857

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 646
$ws.Range("I16").Value = 569
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 569
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -282
$ws.Range("N16").Value = -1374
$ws.Range("H19").Value = 2970.5
$ws.Range("I19").Value = 460.66666
$ws.Range("J19").Value = 10500
$ws.Range("K19").Value = 460.66666
$ws.Range("L19").Value = 10500
$ws.Range("M19").Value = -290.66666
$ws.Range("N19").Value = -10840
$ws.Range("H24").Value = 2970.5
$ws.Range("I24").Value = 460.66666
$ws.Range("J24").Value = 10500
$ws.Range("K24").Value = 460.66666
$ws.Range("L24").Value = 10500
$ws.Range("M24").Value = -290.66666
$ws.Range("N24").Value = -10840
$ws.Range("H86").Value = 142864080
$ws.Range("I86").Value = 200002720
$ws.Range("K86").Value = 200002720
$ws.Range("M86").Value = -200001597
$ws.Range("H89").Value = 142864080
$ws.Range("I89").Value = 200002720
$ws.Range("K89").Value = 1000013600
$ws.Range("M89").Value = -1000007984
$ws.Range("H113").Value = 646
$ws.Range("I113").Value = 569
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 569
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1601
$ws.Range("N113").Value = -5140
$ws.Range("H132").Value = 2114.5557
$ws.Range("I132").Value = 1752.125
$ws.Range("K132").Value = 5256.375
$ws.Range("M132").Value = -2726.375

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H33").Value = 55
$ws.Range("I33").Value = 62.5
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 375
$ws.Range("L33").Value = 240
$ws.Range("M33").Value = -92
$ws.Range("N33").Value = -806
$ws.Range("H81").Value = 183
$ws.Range("I81").Value = 183
$ws.Range("K81").Value = 549
$ws.Range("M81").Value = 574
$ws.Range("H84").Value = 183
$ws.Range("I84").Value = 183
$ws.Range("K84").Value = 1647
$ws.Range("M84").Value = 3969
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242
$ws.Range("H131").Value = 5000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 15000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -25080

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 500000000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 500000000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 500000000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -500000540
$ws.Range("H73").Value = 500000000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 500000000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 500000000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -500001872
$ws.Range("H122").Value = 4474.0454
$ws.Range("I122").Value = 3031.0667
$ws.Range("K122").Value = 9093.2001
$ws.Range("M122").Value = -6643.2001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H23").Value = 999.5
$ws.Range("J23").Value = 999.5
$ws.Range("L23").Value = 999.5
$ws.Range("N23").Value = -1459.5
$ws.Range("H40").Value = 1966.6666
$ws.Range("I40").Value = 1966.6666
$ws.Range("K40").Value = 1966.6666
$ws.Range("M40").Value = -1830.6666
$ws.Range("H51").Value = 33076
$ws.Range("I51").Value = 33076
$ws.Range("K51").Value = 33076
$ws.Range("M51").Value = -32598

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4765
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4710
$ws.Range("H51").Value = 9742.25
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16020
$ws.Range("H52").Value = 12568.4
$ws.Range("I52").Value = 11960.5
$ws.Range("K52").Value = 11960.5
$ws.Range("M52").Value = -11734.5
$ws.Range("H96").Value = 1857
$ws.Range("J96").Value = 2666.3333
$ws.Range("L96").Value = 2666.3333
$ws.Range("N96").Value = -5412.3333
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -39178
